# Bump the published term version and update the publication date
# on the "Metadata" sheet of the ValueSet workbook.
#   Version: 1.0.0 -> 1.1.0
#   Date:    2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
